$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. At the end of the document, split the old "Prompt: Create a cartoon..."
#    paragraph into two paragraphs:
#      a) a bold paragraph with the page title text
#      b) an italic paragraph with the (shortened) meta-description text
#    which used to be the old meta-description content.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)

# 2a. Append the new italic paragraph right after the (still intact) prompt
#     paragraph -- inserting at the absolute end of the document content is
#     a clean split that does not disturb the existing paragraph.
$newText = "Read our review and play Gaelic Gold for free, powered by Nolimit City, featuring a lucky Irish theme and expanding wilds with increasing multipliers."
$appendRange = $d.Range($promptPara.Range.End, $promptPara.Range.End)
$appendXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>" + $newText + "</w:t></w:r></w:p>"
$appendRange.InsertXML($appendXml)

# 2b. Replace the original prompt paragraph's content (title text, bold)
#     while reusing the paragraph's own pre-existing leading empty run.
$promptPara2 = $d.Paragraphs($count)
$boldText = "Play Gaelic Gold for Free - Review of Nolimit City's Slot Game"
$replaceRange = $d.Range($promptPara2.Range.Start, $promptPara2.Range.End)
$replaceXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:b/></w:rPr><w:t>" + $boldText + "</w:t></w:r></w:p>"
$replaceRange.InsertXML($replaceXml)
